$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking text values (prices, percentages) are stored as text,
# matching the original inline-string cell format used throughout the sheet.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "317.83"
$ws.Range("E2").Value = "3.96%"
$ws.Range("D3").Value = "39.70"
$ws.Range("E3").Value = "2.01%"
$ws.Range("D4").Value = "5.132"
$ws.Range("E4").Value = "0.65%"
$ws.Range("D5").Value = "0.08202"
$ws.Range("E5").Value = "1.49%"
$ws.Range("D6").Value = "2.052"
$ws.Range("E6").Value = "6.64%"
$ws.Range("D7").Value = "8.349"
$ws.Range("E7").Value = "3.83%"
$ws.Range("B8").Value = "GateToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D8").Value = "4.312"
$ws.Range("E8").Value = "2.62%"
$ws.Range("B9").Value = "MXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D9").Value = "0.9403"
$ws.Range("E9").Value = "1.39%"
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10").Value = "0.1357"
$ws.Range("E10").Value = "-6.98%"
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").Value = "0.1989"
$ws.Range("E11").Value = "4.12%"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "0.09145"
$ws.Range("E12").Value = "1.32%"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "0.03497"
$ws.Range("E13").Value = "-0.32%"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "0.09788"
$ws.Range("E14").Value = "0.25%"
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").Value = "0.001405"
$ws.Range("E15").Value = "0.27%"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").Value = "0.006011"
$ws.Range("E16").Value = "2.45%"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").Value = "3.688"
$ws.Range("E17").Value = "-2.54%"
$ws.Range("D18").Value = "3.232"
$ws.Range("E18").Value = "-5.18%"
$ws.Range("D19").Value = "0.3481"
$ws.Range("E19").Value = "0.55%"
$ws.Range("D20").Value = "0.1321"
$ws.Range("E20").Value = "-0.47%"
$ws.Range("D21").Value = "4.989"
$ws.Range("E21").Value = "6.32%"
$ws.Range("E22").Value = "1.26%"
$ws.Range("D23").Value = "0.04348"
$ws.Range("E23").Value = "-0.67%"
$ws.Range("E24").Value = "-0.40%"
$ws.Range("E25").Value = "12.47%"
$ws.Range("E26").Value = "-0.12%"
$ws.Range("D27").Value = "0.0003998"
$ws.Range("E27").Value = "-10.12%"
$ws.Range("D39").Value = "0.02257"
$ws.Range("E39").Value = "11.11%"
$ws.Range("D40").Value = "0.05213"
$ws.Range("E40").Value = "3.10%"
$ws.Range("D41").Value = "0.007756"
$ws.Range("E41").Value = "3.10%"
$ws.Range("D42").Value = "0.009861"
$ws.Range("E42").Value = "-0.51%"
$ws.Range("D43").Value = "0.1406"
$ws.Range("E43").Value = "4.74%"
$ws.Range("E44").Value = "-2.88%"
$ws.Range("D45").Value = "0.009333"
$ws.Range("E45").Value = "-6.01%"
$ws.Range("D46").Value = "0.00006609"
$ws.Range("E46").Value = "6.57%"
$ws.Range("D47").Value = "0.00000000750"
$ws.Range("E47").Value = "-0.23%"
$ws.Range("B48").Value = "CoinbaseStockToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
$ws.Range("D48").Value = "0.001689"
$ws.Range("E48").Value = "-6.33%"
$ws.Range("B49").Value = "BOLO"
$ws.Range("C49").Value = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
$ws.Range("D49").Value = "0.002945"
$ws.Range("E49").Value = "2.42%"
$ws.Range("D50").Value = "0.00002099"
$ws.Range("E50").Value = "-0.23%"
$ws.Range("D51").Value = "0.0001999"
$ws.Range("E51").Value = "-0.23%"
